# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# All D/E (Price / Volume(1h)) cells are stored as literal text in this sheet,
# so numeric-looking values are written with NumberFormat forced to "@" (and
# the cell Style restored immediately after) to stop Excel's automatic
# number-typing from turning e.g. "9.10" into the number 9.1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.829.95"
$ws.Range("E2").Value = "  -3.39%  "
$ws.Range("D3").Value = "3.436.24"
$ws.Range("E3").Value = "  -3.17%  "
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  +0.02%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.43"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.94%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.74"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -8.44%  "
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  +0.10%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.623"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -2.10%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +4.49%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.73"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("E12").Value = "  -0.08%  "
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.10"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -3.83%  "
$ws.Range("D14").Value = "3.986.58"
$ws.Range("E14").Value = "  -2.83%  "
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "3.434.36"
$ws.Range("E16").Value = "  -3.18%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.04"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -1.60%  "
$ws.Range("D18").Value = "64.774.90"
$ws.Range("E18").Value = "  -3.34%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.82"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -2.13%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.985"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -1.55%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "406.40"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -6.35%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.19"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +0.67%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.41"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +6.39%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.44"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -2.44%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.22"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +7.70%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.78"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -3.43%  "
$ws.Range("E27").Value = "  -4.00%  "
$ws.Range("E28").Value = "  -2.40%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.93"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -2.56%  "
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.81"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("E32").Value = "  -2.45%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "584.00"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -9.11%  "
$ws.Range("E34").Value = "  -3.42%  "
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.51"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("E36").Value = "  +3.57%  "
$ws.Range("E37").Value = "  +0.29%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.54"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +4.58%  "
$ws.Range("D39").Value = "0.0₃0767"
$ws.Range("E39").Value = "  -6.13%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.08"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -6.69%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.375"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -4.41%  "
$ws.Range("D42").Value = "3.177.33"
$ws.Range("E42").Value = "  +4.46%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  +1.09%  "
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.25"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -2.79%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.50"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -6.62%  "
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("E48").Value = "  -1.61%  "
$ws.Range("E49").Value = "  -4.82%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.42"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -2.92%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.84"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -3.67%  "
